$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure columns B-E are treated as plain text so numeric-looking strings
# (e.g. prices like "213.30") are preserved exactly, matching the source data.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.164.09'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +1.04%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.609.66'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.62%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.996'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -1.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '213.30'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +1.97%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.997'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -1.05%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.484'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +1.80%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.251'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +3.17%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0621'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +2.05%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '18.66'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +5.13%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0791'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +0.31%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.813.15'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -0.50%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.590.48'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -0.54%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.07'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +1.04%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.516'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +1.69%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '26.095.83'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +0.86%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0₃0734'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +2.49%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '60.98'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.60%  '
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -0.54%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '199.36'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +5.45%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.28'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +2.83%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.52'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +2.70%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.03'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +2.05%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.132'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +3.53%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '141.99'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +0.21%  '
$ws.Range("B26").NumberFormat = "@"
$ws.Range("B26").Value = 'Toncoin'
$ws.Range("C26").NumberFormat = "@"
$ws.Range("C26").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.73'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +1.83%  '
$ws.Range("B27").NumberFormat = "@"
$ws.Range("B27").Value = 'BinanceUSD'
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.996'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -1.24%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.32'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +2.87%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.53'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +0.70%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.18'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -2.12%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0478'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +2.38%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.17'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +3.63%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.06'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +2.58%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.53'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +5.26%  '
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -3.29%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.111.02'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +0.45%  '
$ws.Range("B37").NumberFormat = "@"
$ws.Range("B37").Value = 'VeChain'
$ws.Range("C37").NumberFormat = "@"
$ws.Range("C37").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0154'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +2.25%  '
$ws.Range("B38").NumberFormat = "@"
$ws.Range("B38").Value = 'MXToken'
$ws.Range("C38").NumberFormat = "@"
$ws.Range("C38").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.36'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +0.26%  '
$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").Value = 'ImmutableX'
$ws.Range("C39").NumberFormat = "@"
$ws.Range("C39").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.511'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +3.64%  '
$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = 'ARBITRUM'
$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.797'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +0.41%  '
$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = 'PaxDollar'
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.998'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -0.94%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.806'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +8.80%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.16'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +1.88%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '93.16'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -2.35%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.727.65'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -0.43%  '
$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = 'RenderToken'
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.58'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +8.23%  '
$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = 'BabyDogeCoin'
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0₆0107'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -1.27%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '53.99'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +2.09%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0510'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +0.06%  '
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -0.74%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.996'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -1.14%  '
